$wb = $excel.ActiveWorkbook

# --- Sheet 1 (链表) view state: change selection, keep it inactive at the end ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F12").Select() | Out-Null

# --- Sheet 2 (哈希): add new row of content (row 6) ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(6,1).Value = 5
$ws2.Cells.Item(6,2).Value = 350
$ws2.Cells.Item(6,3).Value = "给定两个数组，编写一个函数来计算它们的交集。 "
$ws2.Cells.Item(6,4).Value = "1 选择一个数组，计算每个元素出现的次数`n2 遍历另一个数组，这个元素是否在map中`n3 如果不存在，就迭代下一个元素`n4 如果存在，判断这个数字对应的val是否是0`n5 如果val不是0，就将此元素加入列表`n6 如果val是0，说明公共元素已经找完了`n7 返回列表"
$ws2.Cells.Item(6,5).Value = "哈希表`n交集`n重复元素"

$f6 = $ws2.Cells.Item(6,6)
$f6.Value = "O(M+N),M,N是两个数组的元素个数"
$f6.Font.Name = "Calibri"
$f6.Font.Size = 14
$f6.Font.Color = 0
$f6.HorizontalAlignment = -4131
$f6.VerticalAlignment = -4160
$f6.WrapText = $true
$chars = $f6.Characters(11, 10)
$chars.Font.Name = "Microsoft YaHei"
$chars.Font.Size = 14
$chars.Font.Color = 0

$ws2.Cells.Item(6,7).Value = "O(n)，其中一个数组要将其对应的次数存入哈希表，n是这个数组的元素个数"

$ws2.Rows.Item(6).RowHeight = 154

# --- Sheet 2 (哈希) view state: change selection; activate last so it stays the active tab ---
$ws2.Range("D13").Select() | Out-Null
